# Add a new "State" column (column D) with its value to both worksheets
# in the workbook, as described by the commit diff:
#   - sharedStrings gains "State" and "STATE OF MndstcT TEXAS THL"
#   - both sheets' dimension grows from A1:C2 to A1:D2
#   - D1 = "State" (header), D2 = "STATE OF MndstcT TEXAS THL" (value)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("D1").Value = "State"
    $ws.Range("D2").Value = "STATE OF MndstcT TEXAS THL"
}
